# Paper Presentation slide 1 ("Content Placeholder 2") updates:
#  1. Turn on "Shrink text on overflow" autofit for the placeholder text box
#     (PowerPoint would compute fontScale/lnSpcReduction once the text no
#     longer fits; enabling normAutofit here reflects that change).
#  2. Rework the "Uploaded link: ..." line so it points at the Google Drive
#     upload instead of the old YouTube link, merging/retyping the runs the
#     way PowerPoint does when you select-and-replace part of a line.
#  3. Remove the extra blank centered paragraph that used to sit right
#     after that line.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- 1. Shrink text on overflow -------------------------------------------------
$tf.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

# --- 2. Update the uploaded-link line -------------------------------------------
# Original runs in that paragraph:
#   "Uploaded "              (rPr b=1 dirty=0 smtClean=0)
#   "link: https"             (rPr b=1 smtClean=0)
#   "://youtu.be/M7RxQRvV22I" (rPr b=1 dirty=0)
# Target runs:
#   "Uploaded link"                                                    (keeps run-1 rPr)
#   ": https://drive.google.com/file/d/1lhFvJBDbLc27p8gckzscuh0alHI_mh3O/view?usp=sharing" (keeps run-3 rPr)

# 2a. Extend the first run: "Uploaded " -> "Uploaded link"
$full = $tr.Text
$idx  = $full.IndexOf("Uploaded ")
$run1 = $tr.Characters($idx + 1, "Uploaded ".Length)
$run1.Text = "Uploaded link"

# 2b. Rewrite the third run (the old URL) in place so it keeps its own
#     formatting, turning it into the new address (with the leading colon
#     that used to belong to run two).
$full = $tr.Text
$idx  = $full.IndexOf("://youtu.be/M7RxQRvV22I")
$run3 = $tr.Characters($idx + 1, "://youtu.be/M7RxQRvV22I".Length)
$run3.Text = ": https://drive.google.com/file/d/1lhFvJBDbLc27p8gckzscuh0alHI_mh3O/view?usp=sharing"

# 2c. Delete what's left of the old second run ("link: https"), now
#     redundant since its pieces were folded into runs one and three above.
$full = $tr.Text
$idx  = $full.IndexOf("link: https")
$run2 = $tr.Characters($idx + 1, "link: https".Length)
$run2.Text = ""

# --- 3. Remove the trailing empty centered paragraph ----------------------------
$blank = $tr.Paragraphs(9, 1)
$blank.Delete()
